# Add runmodes for test suites
$wb = $excel.ActiveWorkbook

# 1. Rename the "test_suite" sheet to "testSuite"
$suiteSheet = $wb.Worksheets.Item("test_suite")
$suiteSheet.Name = "testSuite"

# 2. AddCustomerTest: rows 4 and 5 in the runmode column (E) go from "N" to "y"
$addCustomerSheet = $wb.Worksheets.Item("AddCustomerTest")
$addCustomerSheet.Range("E4").Value = "y"
$addCustomerSheet.Range("E5").Value = "y"
$addCustomerSheet.Range("E6").Select()

# 3. OpenAccountTest: add a new "runmode" column (C) with header + value "n"
$openAccountSheet = $wb.Worksheets.Item("OpenAccountTest")
$openAccountSheet.Range("C1").Value = "runmode"
$openAccountSheet.Range("C2").Value = "n"
$openAccountSheet.Range("C2").Select()

# 4. testSuite: runmode column (B) data rows go from "Y" to "y"
$suiteSheet.Range("B2").Value = "y"
$suiteSheet.Range("B3").Value = "y"
$suiteSheet.Range("B4").Value = "y"
$suiteSheet.Range("B5").Select()

# Keep OpenAccountTest as the active/selected tab, matching the original workbook state
$openAccountSheet.Activate()
$openAccountSheet.Range("C2").Select()
